$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "277.41"
Set-TextValue "E2" "1.63%"
Set-TextValue "D3" "27.16"
Set-TextValue "E3" "1.21%"
Set-TextValue "D4" "4.940"
Set-TextValue "E4" "0.81%"
Set-TextValue "D5" "0.06409"
Set-TextValue "E5" "1.47%"
Set-TextValue "E6" "0.49%"
Set-TextValue "D7" "1.251"
Set-TextValue "E7" "-11.26%"
Set-TextValue "E8" "-0.14%"
Set-TextValue "E9" "4.53%"
Set-TextValue "D10" "0.05042"
Set-TextValue "E10" "-1.14%"
Set-TextValue "D11" "0.07515"
Set-TextValue "E11" "1.43%"
Set-TextValue "D12" "0.02893"
Set-TextValue "E12" "-8.36%"
Set-TextValue "D13" "0.09013"
Set-TextValue "E13" "-0.27%"
Set-TextValue "D14" "0.001573"
Set-TextValue "E14" "0.29%"
Set-TextValue "D15" "0.0006398"
Set-TextValue "E15" "1.57%"
Set-TextValue "D16" "0.006049"
Set-TextValue "E16" "0.55%"
Set-TextValue "D17" "3.456"
Set-TextValue "E17" "-0.35%"
Set-TextValue "D18" "3.323"
Set-TextValue "E18" "-0.88%"
Set-TextValue "E19" "0.56%"
Set-TextValue "D21" "0.1337"
Set-TextValue "E21" "0.30%"
Set-TextValue "D22" "3.913"
Set-TextValue "E22" "0.19%"
Set-TextValue "D23" "0.04419"
Set-TextValue "E23" "1.79%"
Set-TextValue "D24" "0.001174"
Set-TextValue "D25" "0.003878"
Set-TextValue "E25" "6.18%"
Set-TextValue "E26" "0.14%"
Set-TextValue "E27" "14.05%"
Set-TextValue "D40" "0.04138"
Set-TextValue "E40" "2.33%"
Set-TextValue "D41" "0.006820"
Set-TextValue "E41" "3.31%"
Set-TextValue "E42" "1.34%"
Set-TextValue "D43" "0.002142"
Set-TextValue "E43" "0.62%"
Set-TextValue "E44" "-10.94%"
Set-TextValue "E45" "-2.52%"
Set-TextValue "D46" "1.487"
Set-TextValue "E46" "-36.89%"

Write-Host "Updated cryptos price/volume values"